$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "万向钱潮"
$ws.Range("B2").Value = "山子高科"
$ws.Range("C2").Value = "山子高科"
$ws.Range("A3").Value = "江波龙"
$ws.Range("B3").Value = "万向钱潮"
$ws.Range("C3").Value = "张江高科"
$ws.Range("A4").Value = "三花智控"
$ws.Range("B4").Value = "三花智控"
$ws.Range("C4").Value = "万向钱潮"
$ws.Range("A5").Value = "华友钴业"
$ws.Range("B5").Value = "领益智造"
$ws.Range("C5").Value = "赛力斯"
$ws.Range("A6").Value = "山子高科"
$ws.Range("B6").Value = "赛力斯"
$ws.Range("C6").Value = "天赐材料"
$ws.Range("A7").Value = "天际股份"
$ws.Range("B7").Value = "华友钴业"
$ws.Range("C7").Value = "中电鑫龙"
$ws.Range("A8").Value = "天赐材料"
$ws.Range("B8").Value = "天赐材料"
$ws.Range("C8").Value = "三花智控"
$ws.Range("A9").Value = "领益智造"
$ws.Range("B9").Value = "国盛金控"
$ws.Range("C9").Value = "蓝丰生化"
$ws.Range("A10").Value = "赛力斯"
$ws.Range("B10").Value = "多氟多"
$ws.Range("C10").Value = "江波龙"
$ws.Range("A11").Value = "张江高科"
$ws.Range("B11").Value = "江波龙"
$ws.Range("C11").Value = "深科技"
$ws.Range("A12").Value = "国盛金控"
$ws.Range("B12").Value = "张江高科"
$ws.Range("C12").Value = "方正科技"
$ws.Range("A13").Value = "中电鑫龙"
$ws.Range("B13").Value = "富临精工"
$ws.Range("C13").Value = "天际股份"
$ws.Range("A14").Value = "德明利"
$ws.Range("B14").Value = "天际股份"
$ws.Range("C14").Value = "华建集团"
$ws.Range("A15").Value = "兆易创新"
$ws.Range("B15").Value = "中电鑫龙"
$ws.Range("C15").Value = "国盛金控"
$ws.Range("A16").Value = "多氟多"
$ws.Range("B16").Value = "东方财富"
$ws.Range("C16").Value = "华友钴业"
$ws.Range("A17").Value = "深科技"
$ws.Range("B17").Value = "长电科技"
$ws.Range("C17").Value = "歌尔股份"
$ws.Range("A18").Value = "富临精工"
$ws.Range("B18").Value = "深科技"
$ws.Range("C18").Value = "吉鑫科技"
$ws.Range("A19").Value = "长电科技"
$ws.Range("B19").Value = "兆易创新"
$ws.Range("C19").Value = "领益智造"
$ws.Range("A20").Value = "寒武纪-U"
$ws.Range("B20").Value = "通富微电"
$ws.Range("C20").Value = "富临精工"
$ws.Range("A21").Value = "香农芯创"
$ws.Range("B21").Value = "江西铜业"
$ws.Range("C21").Value = "长电科技"